$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 (the "default" header/footer -> header2.xml / footer2.xml)
# wdHeaderFooterFirstPage = 2 (the "first page" header/footer -> header1.xml / footer1.xml)

# Headers: both contain the BTec logo, renamed image1.jpg -> image2.jpg
$headerPrimary = $sec.Headers.Item(1)
$headerPrimaryPic = $headerPrimary.Range.InlineShapes.Item(1)
$headerPrimaryPic.Name = "image2.jpg"

$headerFirst = $sec.Headers.Item(2)
$headerFirstPic = $headerFirst.Range.InlineShapes.Item(1)
$headerFirstPic.Name = "image2.jpg"

# Footers: both contain the Pearson logo, renamed image2.png -> image1.png
$footerPrimary = $sec.Footers.Item(1)
$footerPrimaryPic = $footerPrimary.Range.InlineShapes.Item(1)
$footerPrimaryPic.Name = "image1.png"

$footerFirst = $sec.Footers.Item(2)
$footerFirstPic = $footerFirst.Range.InlineShapes.Item(1)
$footerFirstPic.Name = "image1.png"
